# Lab meeting calendar update (8/19)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 10: Solomon -> Schuyler, add note ---
$ws.Range("B10").Value = "Schuyler"
$ws.Range("G10").Value = "Schuyler Presentation and lab cleaning"

# --- Row 13: move Justin out of C13, add Joseph in B13 ---
$ws.Range("C13").ClearContents()
$ws.Range("B13").Value = "Joseph"

# --- Row 14: move Guofu out of B14, add Justin in C14, change E14 ---
$ws.Range("B14").ClearContents()
$ws.Range("C14").Value = "Justin"
$ws.Range("E14").Value = "Someone"

# --- Row 13: add note (after E14 so new shared-string ordering matches) ---
$ws.Range("G13").Value = "Rotation student presentations"

# --- Row 15: swap B15/E15 (Justin <-> Guofu) ---
$ws.Range("B15").Value = "Guofu"
$ws.Range("E15").Value = "Justin"

# --- Row 16: Jordan -> Justin (B16), Justin -> Guofu (E16) ---
$ws.Range("B16").Value = "Justin"
$ws.Range("E16").Value = "Guofu"

# --- New rows 17-20: copy formatting from row 16 down through row 20 ---
$ws.Range("A16").Copy()
$ws.Range("A17:A20").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 17
$ws.Range("A17").Value = 45943
$ws.Range("B17").Value = "Daniel"
$ws.Range("E17").Value = "Justin"
$ws.Range("F17").Formula = "=FALSE"
$ws.Range("G17").Value = "Practice qualifying exam"

# Row 18
$ws.Range("A18").Value = 45950
$ws.Range("B18").Value = "Solomon"
$ws.Range("E18").Value = "Daniel"
$ws.Range("F18").Formula = "=FALSE"

# Row 19
$ws.Range("A19").Value = 45957

# Row 20
$ws.Range("A20").Value = 45964

# --- View state: zoom to 130% and move selection to N7 ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("N7").Select()
